$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2625.5454
$ws.Range("I43").Value = 3979.4
$ws.Range("J43").Value = 1497.3334
$ws.Range("K43").Value = 3979.4
$ws.Range("L43").Value = 1497.3334
$ws.Range("M43").Value = -3910.4
$ws.Range("N43").Value = -1635.3334

$ws.Range("H51").Value = 4999.643
$ws.Range("J51").Value = 4999.643
$ws.Range("L51").Value = 4999.643
$ws.Range("N51").Value = -5967.643

$ws.Range("H98").Value = 414.5
$ws.Range("I98").Value = 414.5
$ws.Range("K98").Value = 414.5
$ws.Range("M98").Value = 1083.5

$ws.Range("H110").Value = 59990.5
$ws.Range("J110").Value = 59990.5
$ws.Range("L110").Value = 59990.5
$ws.Range("N110").Value = -68170.5

$ws.Range("H113").Value = 7375
$ws.Range("I113").Value = 2250
$ws.Range("J113").Value = 12500
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 12500
$ws.Range("M113").Value = 1004
$ws.Range("N113").Value = -19008

$ws.Range("H116").Value = 4876.25
$ws.Range("I116").Value = 4752.5
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4752.5
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1310.5
$ws.Range("N116").Value = -11884

$ws.Range("H122").Value = 414.5
$ws.Range("I122").Value = 414.5
$ws.Range("K122").Value = 1243.5
$ws.Range("M122").Value = 1206.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16955570
$ws.Range("I32").Value = 20411842
$ws.Range("K32").Value = 20411842
$ws.Range("M32").Value = -20411555

$ws.Range("H132").Value = 21741786
$ws.Range("I132").Value = 2703.634
$ws.Range("K132").Value = 8110.902
$ws.Range("M132").Value = -5580.902

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 29990
$ws.Range("J6").Value = 29990
$ws.Range("L6").Value = 29990
$ws.Range("N6").Value = -30216

$ws.Range("H105").Value = 13707.375
$ws.Range("I105").Value = 17576.5
$ws.Range("K105").Value = 17576.5
$ws.Range("M105").Value = -15829.5

$ws.Range("H116").Value = 35999.5
$ws.Range("J116").Value = 35999.5
$ws.Range("L116").Value = 35999.5
$ws.Range("N116").Value = -45177.5

$ws.Range("H141").Value = 20001
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 18492
$ws.Range("J28").Value = 18492
$ws.Range("L28").Value = 18492
$ws.Range("N28").Value = -18982

$ws.Range("H31").Value = 23813672
$ws.Range("I31").Value = 2853.2415
$ws.Range("J31").Value = 76930110
$ws.Range("K31").Value = 2853.2415
$ws.Range("L31").Value = 76930110
$ws.Range("M31").Value = -2558.2415
$ws.Range("N31").Value = -76930700

$ws.Range("H34").Value = 23813672
$ws.Range("I34").Value = 2853.2415
$ws.Range("J34").Value = 76930110
$ws.Range("K34").Value = 2853.2415
$ws.Range("L34").Value = 76930110
$ws.Range("M34").Value = -2651.2415
$ws.Range("N34").Value = -76930514

$ws.Range("H43").Value = 90599.39999999999
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 90599.39999999999
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 90599.39999999999
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -90967.39999999999

$ws.Range("H101").Value = 90599.39999999999
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 90599.39999999999
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 90599.39999999999
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -97089.39999999999

$ws.Range("H107").Value = 1097
$ws.Range("I107").Value = 523.36365
$ws.Range("K107").Value = 523.36365
$ws.Range("M107").Value = 1396.63635

$ws.Range("H141").Value = 397206.28
$ws.Range("I141").Value = 48431.668
$ws.Range("J141").Value = 527996.75
$ws.Range("K141").Value = 48431.668
$ws.Range("L141").Value = 527996.75
$ws.Range("M141").Value = -43251.668
$ws.Range("N141").Value = -538356.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1299.75
$ws.Range("I103").Value = 100.5
$ws.Range("J103").Value = 2499
$ws.Range("K103").Value = 301.5
$ws.Range("L103").Value = 7497
$ws.Range("M103").Value = 577.5
$ws.Range("N103").Value = -9255

$ws.Range("H131").Value = 1619
$ws.Range("I131").Value = 1031.3846
$ws.Range("J131").Value = 1924.56
$ws.Range("K131").Value = 3094.1538
$ws.Range("L131").Value = 5773.68
$ws.Range("M131").Value = 1945.8462
$ws.Range("N131").Value = -15853.68

$ws.Range("H133").Value = 10114.458
$ws.Range("I133").Value = 5102
$ws.Range("J133").Value = 15126.917
$ws.Range("K133").Value = 15306
$ws.Range("L133").Value = 45380.751
$ws.Range("M133").Value = -10246
$ws.Range("N133").Value = -55500.751

$ws.Range("H134").Value = 3794.7778
$ws.Range("I134").Value = 1682.8
$ws.Range("K134").Value = 5048.4
$ws.Range("M134").Value = 21.60000000000036

$ws.Range("H137").Value = 5952.2
$ws.Range("I137").Value = 3612
$ws.Range("J137").Value = 7512.3335
$ws.Range("K137").Value = 10836
$ws.Range("L137").Value = 22537.0005
$ws.Range("M137").Value = -5736
$ws.Range("N137").Value = -32737.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3314.85
$ws.Range("I113").Value = 2484.2307
$ws.Range("J113").Value = 4857.4287
$ws.Range("K113").Value = 2484.2307
$ws.Range("L113").Value = 4857.4287
$ws.Range("M113").Value = -314.2307000000001
$ws.Range("N113").Value = -9197.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 69000
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 4451.722
$ws.Range("I122").Value = 4104.0454
$ws.Range("J122").Value = 4998.0713
$ws.Range("K122").Value = 12312.1362
$ws.Range("L122").Value = 14994.2139
$ws.Range("M122").Value = -9862.136200000001
$ws.Range("N122").Value = -19894.2139

$ws.Range("H136").Value = 2052.362
$ws.Range("I136").Value = 1652.537
$ws.Range("J136").Value = 7450
$ws.Range("K136").Value = 4957.611
$ws.Range("L136").Value = 22350
$ws.Range("M136").Value = -2407.611
$ws.Range("N136").Value = -27450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3944.394
$ws.Range("I132").Value = 3944.394
$ws.Range("K132").Value = 11833.182
$ws.Range("M132").Value = -9303.181999999999
